$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 3
    3  = 7
    4  = 4
    5  = 2
    6  = 8
    7  = 7
    8  = 3
    9  = 5
    10 = 6
    11 = 5
    12 = 1
    13 = 8
    14 = 5
    15 = 6
    16 = 3
    17 = 4
    18 = 4
    19 = 10
    20 = 5
    21 = 8
    22 = 6
    23 = 6
    24 = 9
    25 = 9
    26 = 13
    27 = 10
    28 = 7
    29 = 2
    30 = 10
    31 = 7
    32 = 4
    33 = 3
    34 = 7
    35 = 5
    36 = 3
    37 = 4
    38 = 3
    39 = 1
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
